$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B14").Value = "2"
$ws.Range("B15").Value = "2"
$ws.Range("B16").Value = "2"

$ws.Range("B16").Select()
